$wb = $excel.ActiveWorkbook

# The same set of updates applies to both the "展览" and "全部类型" sheets,
# which mirror each other's data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2995
    $ws.Range("F5").Value = 6741
    $ws.Range("F6").Value = 1739
    $ws.Range("F9").Value = 60
    $ws.Range("F10").Value = 122
    $ws.Range("G10").Value = 55
    $ws.Range("F11").Value = 8
    $ws.Range("F12").Value = 26
}
